$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.222610440090321
$ws.Range("D2").Value = 0.0408356919643964
$ws.Range("E2").Value = 0.1408515772319348
$ws.Range("F2").Value = 0.8187474935530261
$ws.Range("G2").Value = 0.002447370423477938
$ws.Range("I2").Value = 0.840289972293256
$ws.Range("K2").Value = 1.92459884011231
$ws.Range("L2").Value = 0.1707463507152838
$ws.Range("O2").Value = 2.881846131472187

# Row 3
$ws.Range("C3").Value = 0.2156025268874373
$ws.Range("D3").Value = 0.03881409652007761
$ws.Range("E3").Value = 0.1369602421329823
$ws.Range("F3").Value = 0.8259581049857445
$ws.Range("G3").Value = 0.002450456556128578
$ws.Range("I3").Value = 0.8479275664889556
$ws.Range("K3").Value = 1.707805448715021
$ws.Range("L3").Value = 0.1664773415516549
$ws.Range("O3").Value = 2.922367808684896

# Row 4
$ws.Range("C4").Value = 0.2113790143419436
$ws.Range("D4").Value = 0.03756990104601954
$ws.Range("E4").Value = 0.1346409611303372
$ws.Range("F4").Value = 0.8311183509506535
$ws.Range("G4").Value = 0.002452450581892495
$ws.Range("I4").Value = 0.8533021342471798
$ws.Range("K4").Value = 1.574324779179278
$ws.Range("L4").Value = 0.1639583650530909
$ws.Range("O4").Value = 2.949787605068167

# Row 5
$ws.Range("C5").Value = 0.2096779570990037
$ws.Range("D5").Value = 0.03706218268145989
$ws.Range("E5").Value = 0.1337134613088296
$ws.Range("F5").Value = 0.8334050360551757
$ws.Range("G5").Value = 0.002453288165218113
$ws.Range("I5").Value = 0.855664219497001
$ws.Range("K5").Value = 1.519841021617992
$ws.Range("L5").Value = 0.1629575445431684
$ws.Range("O5").Value = 2.961598335440272

# Row 6
$ws.Range("C6").Value = 0.2093967122584672
$ws.Range("D6").Value = 0.0369778353580017
$ws.Range("E6").Value = 0.1335605158512649
$ws.Range("F6").Value = 0.833795828803801
$ws.Range("G6").Value = 0.002453428757686784
$ws.Range("I6").Value = 0.8560668141343584
$ws.Range("K6").Value = 1.510788742630893
$ws.Range("L6").Value = 0.1627929096947582
$ws.Range("O6").Value = 2.963597921155127

# Row 7
$ws.Range("C7").Value = 0.2113559919399677
$ws.Range("D7").Value = 0.0375630565571754
$ws.Range("E7").Value = 0.134628381148552
$ws.Range("F7").Value = 0.8311484462139163
$ws.Range("G7").Value = 0.002452461776603593
$ws.Range("I7").Value = 0.8533332945937531
$ws.Range("K7").Value = 1.573590349205176
$ws.Range("L7").Value = 0.1639447636737188
$ws.Range("O7").Value = 2.949944311884977

# Row 8
$ws.Range("C8").Value = 0.2201776939791955
$ws.Range("D8").Value = 0.04013927522908745
$ws.Range("E8").Value = 0.1394953238958792
$ws.Range("F8").Value = 0.8210813269024726
$ws.Range("G8").Value = 0.002448413993396915
$ws.Range("I8").Value = 0.8427810165201919
$ws.Range("K8").Value = 1.849926951514476
$ws.Range("L8").Value = 0.1692531686931318
$ws.Range("O8").Value = 2.895290013755741

# Row 9
$ws.Range("C9").Value = 0.2381035733915553
$ws.Range("D9").Value = 0.04516662199195309
$ws.Range("E9").Value = 0.1495947009675831
$ws.Range("F9").Value = 0.8071757962825572
$ws.Range("G9").Value = 0.002441259387147193
$ws.Range("I9").Value = 0.8275401473166895
$ws.Range("K9").Value = 2.388775598477139
$ws.Range("L9").Value = 0.1804756024921801
$ws.Range("O9").Value = 2.808333105294423

# Row 10
$ws.Range("C10").Value = 0.2516529295131704
$ws.Range("D10").Value = 0.04884375436763122
$ws.Range("E10").Value = 0.1573539597413358
$ws.Range("F10").Value = 0.8005467174282472
$ws.Range("G10").Value = 0.002436475432466468
$ws.Range("I10").Value = 0.8196898980291749
$ws.Range("K10").Value = 2.782688081240053
$ws.Range("L10").Value = 0.1892195825153209
$ws.Range("O10").Value = 2.756871347124559

# Row 11
$ws.Range("C11").Value = 0.2578987604857446
$ws.Range("D11").Value = 0.05051272951696717
$ws.Range("E11").Value = 0.1609577334195862
$ws.Range("F11").Value = 0.7983163053947351
$ws.Range("G11").Value = 0.002434400659000486
$ws.Range("I11").Value = 0.8168505179802281
$ws.Range("K11").Value = 2.961436533478548
$ws.Range("L11").Value = 0.1933065813619805
$ws.Range("O11").Value = 2.736179386900801

# Row 12
$ws.Range("C12").Value = 0.2602756366570418
$ws.Range("D12").Value = 0.05114415470824696
$ws.Range("E12").Value = 0.1623330342080109
$ws.Range("F12").Value = 0.7975851211617098
$ws.Range("G12").Value = 0.002433629509623773
$ws.Range("I12").Value = 0.8158809466882033
$ws.Range("K12").Value = 3.029057326348379
$ws.Range("L12").Value = 0.1948699836909356
$ws.Range("O12").Value = 2.728736513344472

# Row 13
$ws.Range("C13").Value = 0.2597632143604187
$ws.Range("D13").Value = 0.05100819235387632
$ws.Range("E13").Value = 0.1620363660919324
$ws.Range("F13").Value = 0.7977375419996733
$ws.Range("G13").Value = 0.002433794945651636
$ws.Range("I13").Value = 0.8160850561945665
$ws.Range("K13").Value = 3.014497029652375
$ws.Range("L13").Value = 0.1945325763861661
$ws.Range("O13").Value = 2.730321973587962

# Row 14
$ws.Range("C14").Value = 0.2580940732857755
$ws.Range("D14").Value = 0.05056468907172018
$ws.Range("E14").Value = 0.1610706673062623
$ws.Range("F14").Value = 0.7982538742882781
$ws.Range("G14").Value = 0.002434336925297492
$ws.Range("I14").Value = 0.8167686311895608
$ws.Range("K14").Value = 2.967001106657165
$ws.Range("L14").Value = 0.1934348877255729
$ws.Range("O14").Value = 2.735559175927307

# Row 15
$ws.Range("C15").Value = 0.2570731989688966
$ws.Range("D15").Value = 0.05029295394879796
$ws.Range("E15").Value = 0.1604805327136702
$ws.Range("F15").Value = 0.7985849291932894
$ws.Range("G15").Value = 0.002434670792745594
$ws.Range("I15").Value = 0.8172011104716006
$ws.Range("K15").Value = 2.937899607013208
$ws.Range("L15").Value = 0.1927645728710417
$ws.Range("O15").Value = 2.73881831227942

# Row 16
$ws.Range("C16").Value = 0.2512463945623438
$ws.Range("D16").Value = 0.04873460383858941
$ws.Range("E16").Value = 0.1571199331342257
$ws.Range("F16").Value = 0.8007083241826862
$ws.Range("G16").Value = 0.002436613059743996
$ws.Range("I16").Value = 0.8198902183166155
$ws.Range("K16").Value = 2.770997210586188
$ws.Range("L16").Value = 0.1889546894009442
$ws.Range("O16").Value = 2.758278471026898

# Row 17
$ws.Range("C17").Value = 0.2476928113103156
$ws.Range("D17").Value = 0.047777614002527
$ws.Range("E17").Value = 0.1550772665953559
$ws.Range("F17").Value = 0.80221243101645
$ws.Range("G17").Value = 0.002437830517516997
$ws.Range("I17").Value = 0.8217276116324541
$ws.Range("K17").Value = 2.66849170855528
$ws.Range("L17").Value = 0.1866454658576941
$ws.Range("O17").Value = 2.770914168584625

# Row 18
$ws.Range("C18").Value = 0.2456566234460809
$ws.Range("D18").Value = 0.04722682601402539
$ws.Range("E18").Value = 0.1539093496309079
$ws.Range("F18").Value = 0.8031514284955961
$ws.Range("G18").Value = 0.00243854032218478
$ws.Range("I18").Value = 0.822853281380695
$ws.Range("K18").Value = 2.609491638697932
$ws.Range("L18").Value = 0.1853275496325182
$ws.Range("O18").Value = 2.778437632373794

# Row 19
$ws.Range("C19").Value = 0.2449685374790818
$ws.Range("D19").Value = 0.04704027922594634
$ws.Range("E19").Value = 0.1535151108135437
$ws.Range("F19").Value = 0.8034820292962692
$ws.Range("G19").Value = 0.002438782293014777
$ws.Range("I19").Value = 0.8232462260964226
$ws.Range("K19").Value = 2.589508187896058
$ws.Range("L19").Value = 0.1848830920830693
$ws.Range("O19").Value = 2.781028816219163

# Row 20
$ws.Range("C20").Value = 0.24807029605509
$ws.Range("D20").Value = 0.04787952398459083
$ws.Range("E20").Value = 0.155293990616677
$ws.Range("F20").Value = 0.802044666751101
$ws.Range("G20").Value = 0.002437699928636728
$ws.Range("I20").Value = 0.821524889262605
$ws.Range("K20").Value = 2.679407926702652
$ws.Range("L20").Value = 0.1868902216911579
$ws.Range("O20").Value = 2.769542595027588

# Row 21
$ws.Range("C21").Value = 0.258584023309453
$ws.Range("D21").Value = 0.05069497275645318
$ws.Range("E21").Value = 0.1613540280980956
$ws.Range("F21").Value = 0.7980991325747624
$ws.Range("G21").Value = 0.002434177339214104
$ws.Range("I21").Value = 0.8165649782733482
$ws.Range("K21").Value = 2.980953666312871
$ws.Range("L21").Value = 0.1937568779922714
$ws.Range("O21").Value = 2.734010209602985

# Row 22
$ws.Range("C22").Value = 0.2655235937091049
$ws.Range("D22").Value = 0.05253163687528684
$ws.Range("E22").Value = 0.1653765640384322
$ws.Range("F22").Value = 0.796181860194352
$ws.Range("G22").Value = 0.002431959738922775
$ws.Range("I22").Value = 0.8139393354719715
$ws.Range("K22").Value = 3.177636997827847
$ws.Range("L22").Value = 0.1983364178909142
$ws.Range("O22").Value = 2.713077683091569

# Row 23
$ws.Range("C23").Value = 0.2618136027147102
$ws.Range("D23").Value = 0.0515516973770147
$ws.Range("E23").Value = 0.1632239982450585
$ws.Range("F23").Value = 0.797144460005093
$ws.Range("G23").Value = 0.002433135594305732
$ws.Range("I23").Value = 0.8152841936702941
$ws.Range("K23").Value = 3.072700610345805
$ws.Range("L23").Value = 0.1958838248606156
$ws.Range("O23").Value = 2.724039628839932

# Row 24
$ws.Range("C24").Value = 0.2478996139705885
$ws.Range("D24").Value = 0.04783345235303216
$ws.Range("E24").Value = 0.1551959896290001
$ws.Range("F24").Value = 0.8021202816756272
$ws.Range("G24").Value = 0.002437758937056345
$ws.Range("I24").Value = 0.8216163241218837
$ws.Range("K24").Value = 2.674472916943785
$ws.Range("L24").Value = 0.1867795374077303
$ws.Range("O24").Value = 2.770161876682408

# Row 25
$ws.Range("C25").Value = 0.2331873600559504
$ws.Range("D25").Value = 0.04380937949560604
$ws.Range("E25").Value = 0.1468030295141034
$ws.Range("F25").Value = 0.8103097663709136
$ws.Range("G25").Value = 0.00244311156662566
$ws.Range("I25").Value = 0.8310770920065949
$ws.Range("K25").Value = 2.243341766477329
$ws.Range("L25").Value = 0.1773522813775088
$ws.Range("O25").Value = 2.829682338527988
